# Updates crypto price/volume data per the upstream GitHub Actions refresh,
# and corrects the Kaspa/Mantle row ordering (rows 37-38 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.994.46"
$ws.Range("E2").Value = "  -0.65%  "

$ws.Range("D3").Value = "3.858.06"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'599.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("D6").Value = "'167.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.56%  "

$ws.Range("D7").Value = "3.857.08"
$ws.Range("E7").Value = "  -0.90%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("E10").Value = "  +0.21%  "

$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("E13").Value = "  +1.73%  "

$ws.Range("D14").Value = "'36.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("D15").Value = "4.507.71"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").Value = "3.844.87"
$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("D17").Value = "68.014.59"
$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("E18").Value = "  +7.28%  "

$ws.Range("D19").Value = "'7.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("D21").Value = "'10.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").Value = "'466.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.27%  "

$ws.Range("D23").Value = "'0.730"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.20%  "

$ws.Range("E24").Value = "  -2.97%  "

$ws.Range("D25").Value = "'83.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.08%  "

$ws.Range("D26").Value = "'2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.24%  "

$ws.Range("D27").Value = "'12.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.31%  "

$ws.Range("E28").Value = "  -0.46%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").Value = "'2.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.54%  "

$ws.Range("D31").Value = "4.010.27"
$ws.Range("E31").Value = "  -1.07%  "

$ws.Range("D32").Value = "'7.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "

$ws.Range("E33").Value = "  -1.05%  "

$ws.Range("D34").Value = "'31.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.35%  "

$ws.Range("D35").Value = "3.836.10"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("E36").Value = "  -1.87%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.140"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.51%  "

$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").Value = "'1.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.55%  "

$ws.Range("D39").Value = "'5.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.99%  "

$ws.Range("D40").Value = "'3.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.03%  "

$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").Value = "'0.314"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.78%  "

$ws.Range("D43").Value = "'428.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("E44").Value = "  +0.43%  "

$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").Value = "'47.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.53%  "

$ws.Range("E47").Value = "  +1.60%  "

$ws.Range("D48").Value = "'26.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.02%  "

$ws.Range("D49").Value = "'142.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.94%  "

$ws.Range("E50").Value = "  +6.62%  "

$ws.Range("D51").Value = "'40.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.27%  "
